# "Generate Report for Handback" — refresh the localization-status report:
#   - Overview / per-locale "Status" flips from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - per-locale "Latest Handback DateTime" advances to the new handback run
#   - per-locale "Error Detail" is cleared (no longer stale vs. latest)
#   - a few report columns are widened/narrowed to fit the new content

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText

# Columns E and F grow to fit the longer status text.
$ov.Columns.Item(5).ColumnWidth = 29.166666666666664
$ov.Columns.Item(6).ColumnWidth = 29.166666666666664

# ---- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("K2").Value = "2016-08-26 16:50:05"
$zh.Range("P2").Value = ""

$zh.Columns.Item(3).ColumnWidth = 29.166666666666664
$zh.Columns.Item(16).ColumnWidth = 12.833333333333332

# ---- de-de sheet ------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("K2").Value = "2016-08-26 16:50:20"
$de.Range("P2").Value = ""

$de.Columns.Item(3).ColumnWidth = 29.166666666666664
$de.Columns.Item(16).ColumnWidth = 12.833333333333332
